# Apply the edits described by the commit "sha1 && reviewed directory structure".
#
# 1) Flip the "done" markers on steps 3 and 4 of the client algorithm from
#    "-" (not done) to "+" (done), matching step 1 which is already "+".
# 2) Add a blank paragraph after the "ctrl-c ... safe exit" server bullet,
#    mirroring the blank-line spacing used elsewhere between sections.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "- 3. Считается хеш сумма (SHA-1) всего пакета данных, и дописывается в конце пакета данных.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "+ 3. Считается хеш сумма (SHA-1) всего пакета данных, и дописывается в конце пакета данных.",
    2) | Out-Null

$d.Content.Find.Execute(
    "- 4. Пакет данных шифруется симметричным алгоритмом, сгенерированным ранее ключом.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "+ 4. Пакет данных шифруется симметричным алгоритмом, сгенерированным ранее ключом.",
    2) | Out-Null

$rng = $d.Content
$rng.Find.Execute(
    "При нажатии ctrl-c - сохранение всех данных из linked list в файл и безопасный выход из программы.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertParagraphAfter()

# Re-anchor a zero-length range just inside the freshly inserted paragraph
# (collapsing right at $rng's position is ambiguous -- it sits on the
# boundary between the previous paragraph and the new one, and would
# resolve back to the previous paragraph) and give it an explicit empty
# run of text, like the other blank-line paragraphs elsewhere in the
# document. It already inherits the surrounding Times/14pt formatting.
$newRng = $d.Range($rng.Start + 1, $rng.Start + 1)
$newPara = $newRng.Paragraphs(1)
$newPara.Range.Text = ""
